{"js": "// The document contains a single table of two-digit \u00f7 one-digit practice\n// problems. Every 4th row holds the 5 problem cells for that \"page\"; the\n// rows in between are blank spacer rows. We replace each problem cell's\n// text in place, addressed by (row, column) so the edit is unambiguous\n// even though some new values collide with other (already-processed or\n// not-yet-processed) old values elsewhere in the table.\nconst replacements = [\n  { row: 0, col: 0, oldText: \"75\u00f79=\", newText: \"23\u00f77=\" },\n  { row: 0, col: 1, oldText: \"75\u00f76=\", newText: \"40\u00f78=\" },\n  { row: 0, col: 2, oldText: \"28\u00f79=\", newText: \"97\u00f78=\" },\n  { row: 0, col: 3, oldText: \"75\u00f77=\", newText: \"82\u00f78=\" },\n  { row: 0, col: 4, oldText: \"66\u00f74=\", newText: \"98\u00f72=\" },\n\n  { row: 4, col: 0, oldText: \"97\u00f75=\", newText: \"91\u00f72=\" },\n  { row: 4, col: 1, oldText: \"49\u00f75=\", newText: \"18\u00f72=\" },\n  { row: 4, col: 2, oldText: \"30\u00f79=\", newText: \"69\u00f77=\" },\n  { row: 4, col: 3, oldText: \"84\u00f76=\", newText: \"85\u00f77=\" },\n  { row: 4, col: 4, oldText: \"74\u00f74=\", newText: \"67\u00f72=\" },\n\n  { row: 8, col: 0, oldText: \"35\u00f74=\", newText: \"37\u00f75=\" },\n  { row: 8, col: 1, oldText: \"14\u00f79=\", newText: \"92\u00f77=\" },\n  { row: 8, col: 2, oldText: \"43\u00f79=\", newText: \"32\u00f77=\" },\n  { row: 8, col: 3, oldText: \"51\u00f73=\", newText: \"86\u00f78=\" },\n  { row: 8, col: 4, oldText: \"69\u00f72=\", newText: \"17\u00f79=\" },\n\n  { row: 12, col: 0, oldText: \"25\u00f79=\", newText: \"83\u00f79=\" },\n  { row: 12, col: 1, oldText: \"67\u00f75=\", newText: \"12\u00f74=\" },\n  { row: 12, col: 2, oldText: \"95\u00f78=\", newText: \"43\u00f74=\" },\n  { row: 12, col: 3, oldText: \"77\u00f75=\", newText: \"52\u00f73=\" },\n  { row: 12, col: 4, oldText: \"44\u00f78=\", newText: \"94\u00f78=\" },\n\n  { row: 16, col: 0, oldText: \"99\u00f77=\", newText: \"84\u00f75=\" },\n  { row: 16, col: 1, oldText: \"75\u00f72=\", newText: \"81\u00f77=\" },\n  { row: 16, col: 2, oldText: \"24\u00f74=\", newText: \"13\u00f78=\" },\n  { row: 16, col: 3, oldText: \"35\u00f72=\", newText: \"75\u00f75=\" },\n  { row: 16, col: 4, oldText: \"32\u00f77=\", newText: \"50\u00f72=\" },\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (const { row, col, oldText, newText } of replacements) {\n  const cell = table.getCell(row, col);\n  cell.load(\"value\");\n  await context.sync();\n\n  const current = (cell.value || \"\").trim();\n  if (current !== oldText) {\n    throw new Error(\n      `Unexpected cell text at (${row}, ${col}): expected \"${oldText}\", found \"${current}\"`\n    );\n  }\n\n  cell.value = newText;\n  await context.sync();\n}\n", "ps1": "# The document contains a single table of two-digit \u00f7 one-digit practice\n# problems. Every 4th row holds the 5 problem cells for that \"page\"; the\n# rows in between are blank spacer rows. We replace each problem cell's\n# text in place, addressed by (row, column) -- Word COM's Table.Cell is\n# 1-based -- so the edit is unambiguous even though some new values\n# collide with other (already-processed or not-yet-processed) old values\n# elsewhere in the table.\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n$replacements = @(\n    @{ Row = 1;  Col = 1; Old = \"75\u00f79=\"; New = \"23\u00f77=\" },\n    @{ Row = 1;  Col = 2; Old = \"75\u00f76=\"; New = \"40\u00f78=\" },\n    @{ Row = 1;  Col = 3; Old = \"28\u00f79=\"; New = \"97\u00f78=\" },\n    @{ Row = 1;  Col = 4; Old = \"75\u00f77=\"; New = \"82\u00f78=\" },\n    @{ Row = 1;  Col = 5; Old = \"66\u00f74=\"; New = \"98\u00f72=\" },\n\n    @{ Row = 5;  Col = 1; Old = \"97\u00f75=\"; New = \"91\u00f72=\" },\n    @{ Row = 5;  Col = 2; Old = \"49\u00f75=\"; New = \"18\u00f72=\" },\n    @{ Row = 5;  Col = 3; Old = \"30\u00f79=\"; New = \"69\u00f77=\" },\n    @{ Row = 5;  Col = 4; Old = \"84\u00f76=\"; New = \"85\u00f77=\" },\n    @{ Row = 5;  Col = 5; Old = \"74\u00f74=\"; New = \"67\u00f72=\" },\n\n    @{ Row = 9;  Col = 1; Old = \"35\u00f74=\"; New = \"37\u00f75=\" },\n    @{ Row = 9;  Col = 2; Old = \"14\u00f79=\"; New = \"92\u00f77=\" },\n    @{ Row = 9;  Col = 3; Old = \"43\u00f79=\"; New = \"32\u00f77=\" },\n    @{ Row = 9;  Col = 4; Old = \"51\u00f73=\"; New = \"86\u00f78=\" },\n    @{ Row = 9;  Col = 5; Old = \"69\u00f72=\"; New = \"17\u00f79=\" },\n\n    @{ Row = 13; Col = 1; Old = \"25\u00f79=\"; New = \"83\u00f79=\" },\n    @{ Row = 13; Col = 2; Old = \"67\u00f75=\"; New = \"12\u00f74=\" },\n    @{ Row = 13; Col = 3; Old = \"95\u00f78=\"; New = \"43\u00f74=\" },\n    @{ Row = 13; Col = 4; Old = \"77\u00f75=\"; New = \"52\u00f73=\" },\n    @{ Row = 13; Col = 5; Old = \"44\u00f78=\"; New = \"94\u00f78=\" },\n\n    @{ Row = 17; Col = 1; Old = \"99\u00f77=\"; New = \"84\u00f75=\" },\n    @{ Row = 17; Col = 2; Old = \"75\u00f72=\"; New = \"81\u00f77=\" },\n    @{ Row = 17; Col = 3; Old = \"24\u00f74=\"; New = \"13\u00f78=\" },\n    @{ Row = 17; Col = 4; Old = \"35\u00f72=\"; New = \"75\u00f75=\" },\n    @{ Row = 17; Col = 5; Old = \"32\u00f77=\"; New = \"50\u00f72=\" }\n)\n\nforeach ($r in $replacements) {\n    $cell = $table.Cell($r.Row, $r.Col)\n    $range = $cell.Range\n    # Cell.Range.Text includes the trailing cell-mark / paragraph-mark\n    # control characters (chr 13 + chr 7) -- strip those before comparing.\n    $current = $range.Text.TrimEnd([char]7, [char]13)\n    if ($current -ne $r.Old) {\n        throw \"Unexpected cell text at row $($r.Row), col $($r.Col): expected '$($r.Old)', found '$current'\"\n    }\n    $range.Text = $r.New\n}\n"}
